$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The "Saldo" export sheet got a new account (FERNANDO) added, an old
# account (MARCELO / 000772433) removed, and the three accounts that
# used to sit right after it (RAFAEL, THEO, NATALIA) moved up to sit
# right after MATEUS's old slot while MATEUS's own balance changed.
# ------------------------------------------------------------------

# 1) Insert a new row for FERNANDO right above LARISSA (row 5).
$ws.Rows.Item(5).Insert()
$newRow = $ws.Cells.Item(5, 1)
$newRow.NumberFormat = "@"
$newRow.Value = "004895776"
$ws.Cells.Item(5, 2).Value = "FERNANDO"
$ws.Cells.Item(5, 3).Value = 31044.77

# 2) Remove the MARCELO (000772433) row, now at row 8.
$ws.Rows.Item(8).Delete()

# 3) Rows 8-11 now hold MATEUS, RAFAEL, THEO, NATALIA (in that order).
#    Rewrite them so RAFAEL/THEO/NATALIA come first, followed by MATEUS
#    with its updated balance.
$accA = $ws.Cells.Item(8, 1)
$accA.NumberFormat = "@"
$accA.Value = "004454365"
$ws.Cells.Item(8, 2).Value = "RAFAEL"
$ws.Cells.Item(8, 3).Value = 13713.49

$accB = $ws.Cells.Item(9, 1)
$accB.NumberFormat = "@"
$accB.Value = "004550750"
$ws.Cells.Item(9, 2).Value = "THEO"
$ws.Cells.Item(9, 3).Value = 10305.28

$accC = $ws.Cells.Item(10, 1)
$accC.NumberFormat = "@"
$accC.Value = "004482102"
$ws.Cells.Item(10, 2).Value = "NATALIA"
$ws.Cells.Item(10, 3).Value = 8500

$accD = $ws.Cells.Item(11, 1)
$accD.NumberFormat = "@"
$accD.Value = "004451652"
$ws.Cells.Item(11, 2).Value = "MATEUS"
$ws.Cells.Item(11, 3).Value = 6006.82
